$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the rows that currently hold the two "Requisitos" entries
# (one starting with "LOM3202", the other with "LOM3206") by scanning
# column B of the used range, instead of assuming fixed row numbers.
$used = $ws.UsedRange
$rowCount = $used.Rows.Count

$row3202 = -1
$row3206 = -1

for ($r = 1; $r -le $rowCount; $r++) {
    $cellVal = $ws.Cells.Item($r, 2).Value2
    if ($cellVal -ne $null) {
        $text = [string]$cellVal
        if ($text.StartsWith("LOM3202")) {
            $row3202 = $r
        }
        elseif ($text.StartsWith("LOM3206")) {
            $row3206 = $r
        }
    }
}

if ($row3202 -gt 0 -and $row3206 -gt 0) {
    $textReq = $ws.Cells.Item($row3202, 2).Value2
    $textInd = $ws.Cells.Item($row3206, 2).Value2

    # Swap the displayed text between the two rows (columns B and C),
    # matching the reordering of the shared-string table entries.
    $ws.Cells.Item($row3202, 2).Value = $textInd
    $ws.Cells.Item($row3202, 3).Value = $textInd
    $ws.Cells.Item($row3206, 2).Value = $textReq
    $ws.Cells.Item($row3206, 3).Value = $textReq
}
